$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency list refresh (GitHub Actions scheduled update).
# Column D holds price text such as "616.87" or "64.867.47"; left alone, Excel
# smart entry would auto-convert the plain-looking values into real numbers.
# Force those specific cells to a Text number format first so the new price
# keeps being stored as text, matching the sheet's original layout.
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "64.867.47"
$ws.Cells.Item(2, 5).Value = "  +1.24%  "
$ws.Cells.Item(3, 4).Value = "3.174.54"
$ws.Cells.Item(3, 5).Value = "  +1.30%  "
$ws.Cells.Item(4, 5).Value = "  +0.10%  "
$ws.Cells.Item(5, 4).Value = "616.87"
$ws.Cells.Item(5, 5).Value = "  +0.88%  "
$ws.Cells.Item(6, 4).Value = "147.04"
$ws.Cells.Item(6, 5).Value = "  -1.65%  "
$ws.Cells.Item(7, 5).Value = "  -0.11%  "
$ws.Cells.Item(8, 4).Value = "3.172.74"
$ws.Cells.Item(8, 5).Value = "  +1.29%  "
$ws.Cells.Item(9, 5).Value = "  -0.27%  "
$ws.Cells.Item(10, 5).Value = "  +0.10%  "
$ws.Cells.Item(11, 5).Value = "  -1.08%  "
$ws.Cells.Item(12, 5).Value = "  -0.46%  "
$ws.Cells.Item(13, 5).Value = "  +0.96%  "
$ws.Cells.Item(14, 4).Value = "35.99"
$ws.Cells.Item(14, 5).Value = "  -2.44%  "
$ws.Cells.Item(15, 4).Value = "3.695.08"
$ws.Cells.Item(15, 5).Value = "  +2.29%  "
$ws.Cells.Item(16, 5).Value = "  +3.17%  "
$ws.Cells.Item(17, 4).Value = "64.868.09"
$ws.Cells.Item(17, 5).Value = "  +1.18%  "
$ws.Cells.Item(18, 4).Value = "3.170.53"
$ws.Cells.Item(18, 5).Value = "  +1.21%  "
$ws.Cells.Item(19, 4).Value = "6.96"
$ws.Cells.Item(19, 5).Value = "  -0.73%  "
$ws.Cells.Item(20, 4).Value = "482.18"
$ws.Cells.Item(20, 5).Value = "  -0.39%  "
$ws.Cells.Item(21, 4).Value = "14.79"
$ws.Cells.Item(22, 5).Value = "  +1.07%  "
$ws.Cells.Item(23, 5).Value = "  +2.59%  "
$ws.Cells.Item(24, 4).Value = "13.87"
$ws.Cells.Item(24, 5).Value = "  +0.25%  "
$ws.Cells.Item(25, 4).Value = "84.80"
$ws.Cells.Item(25, 5).Value = "  +0.52%  "
$ws.Cells.Item(26, 5).Value = "  -0.09%  "
$ws.Cells.Item(27, 5).Value = "  -3.70%  "
$ws.Cells.Item(28, 4).Value = "8.67"
$ws.Cells.Item(28, 5).Value = "  +0.83%  "
$ws.Cells.Item(29, 5).Value = "  -5.31%  "
$ws.Cells.Item(30, 2).Value = "ImmutableX"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(30, 4).Value = "2.12"
$ws.Cells.Item(30, 5).Value = "  -6.16%  "
$ws.Cells.Item(31, 2).Value = "NEARProtocol"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(31, 4).Value = "6.94"
$ws.Cells.Item(31, 5).Value = "  -1.52%  "
$ws.Cells.Item(32, 5).Value = "  +0.30%  "
$ws.Cells.Item(33, 5).Value = "  +0.03%  "
$ws.Cells.Item(34, 4).Value = "26.78"
$ws.Cells.Item(34, 5).Value = "  -0.08%  "
$ws.Cells.Item(35, 5).Value = "  +2.23%  "
$ws.Cells.Item(36, 4).Value = "0.0₃0792"
$ws.Cells.Item(36, 5).Value = "  +5.32%  "
$ws.Cells.Item(37, 4).Value = "6.06"
$ws.Cells.Item(37, 5).Value = "  -1.00%  "
$ws.Cells.Item(38, 5).Value = "  -1.15%  "
$ws.Cells.Item(39, 4).Value = "53.22"
$ws.Cells.Item(39, 5).Value = "  -2.49%  "
$ws.Cells.Item(40, 4).Value = "467.81"
$ws.Cells.Item(40, 5).Value = "  +3.48%  "
$ws.Cells.Item(41, 5).Value = "  +0.25%  "
$ws.Cells.Item(42, 5).Value = "  -3.11%  "
$ws.Cells.Item(43, 4).Value = "8.43"
$ws.Cells.Item(43, 5).Value = "  -1.10%  "
$ws.Cells.Item(44, 4).Value = "2.863.30"
$ws.Cells.Item(44, 5).Value = "  -0.65%  "
$ws.Cells.Item(45, 5).Value = "  +0.40%  "
$ws.Cells.Item(46, 5).Value = "  -1.01%  "
$ws.Cells.Item(47, 5).Value = "  +5.86%  "
$ws.Cells.Item(48, 4).Value = "26.89"
$ws.Cells.Item(48, 5).Value = "  +0.34%  "
$ws.Cells.Item(49, 2).Value = "Arweave"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Cells.Item(49, 4).Value = "36.23"
$ws.Cells.Item(49, 5).Value = "  +8.91%  "
$ws.Cells.Item(50, 2).Value = "USDe"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(50, 4).Value = "1.00"
$ws.Cells.Item(50, 5).Value = "  +0.14%  "
$ws.Cells.Item(51, 2).Value = "Stellar"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(51, 4).Value = "0.115"
$ws.Cells.Item(51, 5).Value = "  -0.81%  "
